# League-base update (24-02-2024 23:13): several match rows in the
# "Germany Bundesliga I" sheet were written against the wrong fixture.
# The fix re-pairs each affected row with the data that actually belongs
# to it — i.e. the row's id/date columns (A, C, D, E) stay put, while the
# rest of the record (B, and F through AC: teams, score, odds, ...) is
# exchanged with the partner row(s) below. The result is a pure swap /
# cyclic permutation of row payloads:
#   293 <-> 294
#   311 -> 315 -> 313 -> 312 -> 314 -> 311  (cycle)
#   325 <-> 328

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Bundesliga I")

# Row -> row whose payload (cols B, F:AC) it should receive.
$mapping = @{
    293 = 294
    294 = 293
    311 = 315
    312 = 314
    313 = 312
    314 = 311
    315 = 313
    325 = 328
    328 = 325
}

# Snapshot every affected row's current payload BEFORE any writes happen,
# so overwriting one row never clobbers data another row still needs.
$original = @{}
foreach ($r in $mapping.Keys) {
    $original[$r] = $ws.Range("B$r`:AC$r").Value2
}

# Now write each row's new payload from the captured snapshot.
foreach ($r in $mapping.Keys) {
    $srcRow = $mapping[$r]
    $ws.Range("B$r`:AC$r").Value = $original[$srcRow]
}
